$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Replace every negative reward value in the B:D range (rows 2-13) with the
# literal text "n" (string), leaving non-negative numbers untouched.
for ($r = 2; $r -le 13; $r++) {
    for ($c = 2; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -lt 0) {
            $cell.Value = "n"
        }
    }
}

# Update the selection shown in the worksheet view.
$ws.Range("D16:D17").Select()
